# json cucumber report added in loaca
# Refresh flight listing data: new (cheaper) Air India flights inserted at the
# top of the sorted-by-price list, remaining rows shift down, and the two
# "Go First" rows are replaced with new IndiGo / SpiceJet fares.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("Air India", "09:55", "03 h 05 m", "₹ 8,248")
    3  = @("Air India", "20:15", "02 h 30 m", "₹ 8,248")
    4  = @("SpiceJet",  "21:40", "02 h 45 m", "₹ 8,982")
    5  = @("SpiceJet",  "08:30", "05 h 30 m", "₹ 8,982")
    6  = @("Air India", "06:00", "02 h 55 m", "₹ 8,983")
    7  = @("IndiGo",    "06:20", "02 h 50 m", "₹ 8,983")
    8  = @("Vistara",   "07:05", "02 h 50 m", "₹ 8,983")
    9  = @("IndiGo",    "08:45", "02 h 45 m", "₹ 8,983")
    10 = @("Vistara",   "10:35", "02 h 50 m", "₹ 8,983")
    11 = @("IndiGo",    "10:40", "02 h 45 m", "₹ 8,983")
    12 = @("IndiGo",    "13:20", "03 h",      "₹ 8,983")
    13 = @("IndiGo",    "15:10", "02 h 50 m", "₹ 8,983")
    14 = @("IndiGo",    "16:35", "02 h 50 m", "₹ 8,983")
    15 = @("Air India", "16:55", "02 h 55 m", "₹ 8,983")
    16 = @("Vistara",   "17:15", "02 h 45 m", "₹ 8,983")
    17 = @("IndiGo",    "18:15", "02 h 50 m", "₹ 8,983")
    18 = @("IndiGo",    "19:35", "02 h 50 m", "₹ 8,983")
    19 = @("IndiGo",    "22:50", "02 h 55 m", "₹ 8,983")
    20 = @("IndiGo",    "11:00", "04 h 35 m", "₹ 9,036")
    21 = @("SpiceJet",  "19:40", "13 h 10 m", "₹ 13,981")
    22 = @("SpiceJet",  "19:40", "13 h 10 m", "₹ 13,509")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
    $ws.Range("D$row").Value = $vals[3]
}
